{"js": "// Update the worksheet date and all \"division problem\" cells to the new\n// values for this day's printout.\nconst replacements = [\n    { find: \"2025-10-11 Saturday\", replace: \"2025-10-12 Sunday\" },\n    { find: \"370\u00f73=\", replace: \"168\u00f76=\" },\n    { find: \"106\u00f74=\", replace: \"227\u00f77=\" },\n    { find: \"750\u00f78=\", replace: \"161\u00f72=\" },\n    { find: \"513\u00f75=\", replace: \"263\u00f73=\" },\n    { find: \"560\u00f78=\", replace: \"318\u00f77=\" },\n    { find: \"499\u00f72=\", replace: \"377\u00f73=\" },\n    { find: \"994\u00f72=\", replace: \"687\u00f74=\" },\n    { find: \"777\u00f76=\", replace: \"785\u00f79=\" },\n    { find: \"809\u00f77=\", replace: \"793\u00f79=\" },\n    { find: \"962\u00f76=\", replace: \"485\u00f72=\" },\n    { find: \"243\u00f76=\", replace: \"756\u00f76=\" },\n    { find: \"631\u00f74=\", replace: \"383\u00f74=\" },\n    { find: \"969\u00f78=\", replace: \"428\u00f77=\" },\n    { find: \"413\u00f72=\", replace: \"114\u00f72=\" },\n    { find: \"260\u00f74=\", replace: \"596\u00f73=\" },\n    { find: \"773\u00f78=\", replace: \"157\u00f77=\" },\n    { find: \"382\u00f72=\", replace: \"482\u00f73=\" },\n    { find: \"283\u00f72=\", replace: \"405\u00f75=\" },\n    { find: \"949\u00f72=\", replace: \"453\u00f72=\" },\n    { find: \"743\u00f79=\", replace: \"870\u00f77=\" },\n    { find: \"743\u00f76=\", replace: \"819\u00f75=\" },\n    { find: \"550\u00f74=\", replace: \"327\u00f74=\" },\n    { find: \"983\u00f73=\", replace: \"570\u00f76=\" },\n    { find: \"753\u00f78=\", replace: \"332\u00f77=\" },\n    { find: \"440\u00f78=\", replace: \"405\u00f74=\" },\n];\n\nfor (const { find, replace } of replacements) {\n    const results = context.document.body.search(find, { matchCase: true, matchWholeWord: false });\n    results.load(\"items\");\n    await context.sync();\n\n    for (const item of results.items) {\n        item.insertText(replace, Word.InsertLocation.replace);\n    }\n    await context.sync();\n}\n", "ps1": "# Update the worksheet date and all \"division problem\" cells to the new\n# values for this day's printout.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{Find = \"2025-10-11 Saturday\"; Replace = \"2025-10-12 Sunday\"},\n    @{Find = \"370\u00f73=\";              Replace = \"168\u00f76=\"},\n    @{Find = \"106\u00f74=\";              Replace = \"227\u00f77=\"},\n    @{Find = \"750\u00f78=\";              Replace = \"161\u00f72=\"},\n    @{Find = \"513\u00f75=\";              Replace = \"263\u00f73=\"},\n    @{Find = \"560\u00f78=\";              Replace = \"318\u00f77=\"},\n    @{Find = \"499\u00f72=\";              Replace = \"377\u00f73=\"},\n    @{Find = \"994\u00f72=\";              Replace = \"687\u00f74=\"},\n    @{Find = \"777\u00f76=\";              Replace = \"785\u00f79=\"},\n    @{Find = \"809\u00f77=\";              Replace = \"793\u00f79=\"},\n    @{Find = \"962\u00f76=\";              Replace = \"485\u00f72=\"},\n    @{Find = \"243\u00f76=\";              Replace = \"756\u00f76=\"},\n    @{Find = \"631\u00f74=\";              Replace = \"383\u00f74=\"},\n    @{Find = \"969\u00f78=\";              Replace = \"428\u00f77=\"},\n    @{Find = \"413\u00f72=\";              Replace = \"114\u00f72=\"},\n    @{Find = \"260\u00f74=\";              Replace = \"596\u00f73=\"},\n    @{Find = \"773\u00f78=\";              Replace = \"157\u00f77=\"},\n    @{Find = \"382\u00f72=\";              Replace = \"482\u00f73=\"},\n    @{Find = \"283\u00f72=\";              Replace = \"405\u00f75=\"},\n    @{Find = \"949\u00f72=\";              Replace = \"453\u00f72=\"},\n    @{Find = \"743\u00f79=\";              Replace = \"870\u00f77=\"},\n    @{Find = \"743\u00f76=\";              Replace = \"819\u00f75=\"},\n    @{Find = \"550\u00f74=\";              Replace = \"327\u00f74=\"},\n    @{Find = \"983\u00f73=\";              Replace = \"570\u00f76=\"},\n    @{Find = \"753\u00f78=\";              Replace = \"332\u00f77=\"},\n    @{Find = \"440\u00f78=\";              Replace = \"405\u00f74=\"}\n)\n\nforeach ($item in $replacements) {\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $item.Find\n    $find.Replacement.Text = $item.Replace\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($null, $null, $null, $null, $null, $null, $true, 1, $null, $null, 2)\n}\n\n$d.Saved = $false\n"}
